$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 82; this shifts the existing rows 82-121
# down to 83-122 (carrying all their values/styles along), matching the
# diff where old row 82 content moves to row 83, old row 83 -> 84, etc.,
# and a brand-new "Primera" Sandia record is inserted as the new row 82.
$ws.Rows.Item(82).Insert()

# Populate the newly inserted row 82 with the new weekly record.
$ws.Cells.Item(82, 1).Value = 11
$ws.Cells.Item(82, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(82, 3).Value = "Bíobío"
$ws.Cells.Item(82, 4).Value = 44904
$ws.Cells.Item(82, 5).Value = 8
$ws.Cells.Item(82, 6).Value = 100112028
$ws.Cells.Item(82, 7).Value = "Sandia"
$ws.Cells.Item(82, 8).Value = "Sin especificar"
$ws.Cells.Item(82, 9).Value = "Primera"
$ws.Cells.Item(82, 10).Value = 2000
$ws.Cells.Item(82, 11).Value = 3000
$ws.Cells.Item(82, 12).Value = 3200
$ws.Cells.Item(82, 13).Value = 3100
$ws.Cells.Item(82, 14).Value = "$/unidad"
$ws.Cells.Item(82, 15).Value = "Paine"
$ws.Cells.Item(82, 16).Value = 3100
$ws.Cells.Item(82, 17).Value = 1
$ws.Cells.Item(82, 18).Value = "Hortaliza"
